$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")
$ws.Activate()

# 1. Delete column I (imageLink) - shifts J..N left to I..M
$ws.Columns.Item(9).Delete()

# 2. Correct bbid value for Dried Cranberry (row 13)
$ws.Range("C13").Value = 1724

# 3. Fill in the "amount_buy" (K) values for every data row (2-19)
$amountBuy = @{
    2 = 22
    3 = 22
    4 = 29
    5 = 24
    6 = 29
    7 = 32
    8 = 20
    9 = 24
    10 = 36
    11 = 24
    12 = 25
    13 = 29
    14 = 31
    15 = 20
    16 = 16
    17 = 34
    18 = 34
    19 = 30
}
foreach ($r in 2..19) {
    $ws.Range("K$r").Value = $amountBuy[$r]
}

# 4. Add the new "order_cost" (M) column formula for rows 2-19
foreach ($r in 2..19) {
    $ws.Range("M$r").Formula = "=H$r*K$r"
}

# 5. Remove the old SUM row (now at row 20, leftover in column L after the
#    column delete shifted the previous M-column SUM formula into L)
$ws.Rows.Item(20).Delete()

# 6. Update the view/selection to match what was left after editing
$ws.Range("C14").Select()
$excel.ActiveWindow.ScrollColumn = 2

Write-Host "Done with main edits"
Write-Host "UsedRange:" $ws.UsedRange.Address()
